$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "nom"
$ws.Range("B1").Value = "prénom"
$ws.Range("D1").Value = "adresse"
